$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 edits
$ws.Range("G7").Clear()
$ws.Range("J7").Value = "Tyhjennysväli ok"
$ws.Range("M7").Value = "pidennetty"
